$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first week start date to the 2024-2025 academic year.
# Dependent cells (E8:AH8) are formulas ("=D8+7", etc.) and recompute automatically.
$ws.Range("D8").Value = 45579

# Fill in the task labels that used to be placeholder ("...") rows further down
# the sheet; they are now given real content at the top of the task list.
# (Written in this order so new shared-string entries land in the same
# index order as the reference workbook: Tampon, Analyse, Rencontres, etc.)
$ws.Range("B10").Value = "Tampon de récupération (dodo)"
$ws.Range("B12").Value = "Analyse"
$ws.Range("B11").Value = "Rencontres"
$ws.Range("B13").Value = "etc"

# Clear out the old placeholder ("...") labels that used to sit on rows 36-39.
$ws.Range("B36").Value = ""
$ws.Range("B37").Value = ""
$ws.Range("B38").Value = ""
$ws.Range("B39").Value = ""
